$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 648 ("「眠りたいです」" entry). All rows below
# shift up by one, matching the diff (dimension shrinks from C842 to C841).
$ws.Rows.Item(648).Delete()
